$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11: A=Colaborador_id, B=Colaborador_nome, C=Departamento,
# D=Motivo_da_ausencia, E=Horas_de_ausencia, F=Data_da_ausencia, G=Salario

$data = @(
    @{ Row = 2;  A = 47332; B = "Raquel das Neves";       C = "Financeiro";  D = "Problemas pessoais";   E = 8; F = 45091; G = 7226.98 }
    @{ Row = 3;  A = 15439; B = "Diogo Aragão";            C = "Jurídico";    D = "Outros";                E = 7; F = 45083; G = 8464.86 }
    @{ Row = 4;  A = 89755; B = "Alana Moura";             C = "Operações";   D = "Doença";                E = 1; F = 45104; G = 4056.14 }
    @{ Row = 5;  A = 97626; B = "Rebeca Cardoso";          C = "Vendas";      D = "Viagem de negócios";   E = 3; F = 45090; G = 5174.3 }
    @{ Row = 6;  A = 19453; B = "Sophie Correia";          C = "Operações";   D = "Viagem de negócios";   E = 1; F = 45080; G = 9112.68 }
    @{ Row = 7;  A = 15576; B = "Lucas Gabriel Cardoso";   C = "TI";          D = "Problemas pessoais";   E = 7; F = 45100; G = 2500.66 }
    @{ Row = 8;  A = 36610; B = "Joaquim da Mata";         C = "Vendas";      D = "Problemas pessoais";   E = 3; F = 45088; G = 3552.61 }
    @{ Row = 9;  A = 40490; B = "Ana Clara Campos";        C = "Vendas";      D = "Outros";                E = 5; F = 45080; G = 6339.72 }
    @{ Row = 10; A = 76302; B = "Enzo Gabriel Costa";      C = "Operações";   D = "Doença";                E = 8; F = 45090; G = 6395.9 }
    @{ Row = 11; A = 51262; B = "Kevin Monteiro";          C = "Engenharia";  D = "Consulta médica";      E = 4; F = 45080; G = 10023.61 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}

$wb.Save()
